$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49 (pushes old rows 49-109 down to 50-110)
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record
$ws.Cells.Item(49, 1).Value  = 11
$ws.Cells.Item(49, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value  = "Bíobío"
$ws.Cells.Item(49, 4).Value  = (Get-Date -Year 2023 -Month 5 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(49, 5).Value  = 8
$ws.Cells.Item(49, 6).Value  = 100112037
$ws.Cells.Item(49, 7).Value  = "Cebollín"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Primera"
$ws.Cells.Item(49, 10).Value = 40
$ws.Cells.Item(49, 11).Value = 4500
$ws.Cells.Item(49, 12).Value = 5000
$ws.Cells.Item(49, 13).Value = 4750
$ws.Cells.Item(49, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 132
$ws.Cells.Item(49, 17).Value = 36
$ws.Cells.Item(49, 18).Value = "Hortaliza"
